$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: free up the shared-string slots that need to be reshuffled ---
# B10:B16 hold the "scheme name" labels from "NoRotation-tilt60deg" through
# "Gaussian-Quadrature" (shared-string indices 8-14). Clearing them drops
# those strings from the table (garbage collected) so we can re-add them
# (plus the 3 new spiral schemes) in the right order.
$ws.Range("B10:B16").ClearContents()

# C2:M2 hold the HKL-tuple column headers (shared-string indices 15-25).
# Clear them too so they get re-appended to the table *after* the scheme
# names below, matching the original build order (all row labels first,
# then the header row).
$ws.Range("C2:M2").ClearContents()

# --- Step 2: rebuild the scheme-name column, in final row order ---
# Rows 10-13 are the newly-run schemes (Gaussian-Quadrature moved up from
# its old spot, plus the three new spiral sweeps); rows 14-16 are the
# schemes that used to sit at rows 10-12.
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# --- Step 3: three brand-new rows (17-19) for the remaining hex-grid schemes ---
$ws.Range("A17:A19").Value = $ws.Range("A16").Value
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C17:M19").Value = 1

# --- Step 4: rebuild the HKL header row (same text, re-appended after scheme names) ---
$ws.Range("C2").Value = "[1, 1, 0]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 1, 1]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[3, 1, 0]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 2, 1]"
$ws.Range("J2").Value = "[4, 0, 0]"
$ws.Range("K2").Value = "2Pairs"
$ws.Range("L2").Value = "4Pairs"
$ws.Range("M2").Value = "MaxUnique"
